# Weekly profitability analysis refresh - Aug 2, 2025
# Rolls the 4-week window forward: a new week ("Jul 26 - Aug 01") is added,
# the oldest week ("Jun 28 - Jul 04") drops off, and all summary metrics
# (tickers/winners/losers/win rate/avg return/best/worst) are recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Period labels (column B) roll forward by one week ---
$ws.Range("B2").Value = "Jul 26 - Aug 01"
$ws.Range("B3").Value = "Jul 19 - Jul 25"
$ws.Range("B4").Value = "Jul 12 - Jul 18"
$ws.Range("B5").Value = "Jul 05 - Jul 11"

# --- Row 2 (Week 4) ---
$ws.Range("C2").Value = 38
$ws.Range("D2").Value = 30
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 78.94736842105263
$ws.Range("G2").Value = 1.693924098935704
$ws.Range("H2").Value = "PEL"
$ws.Range("I2").Value = 7.688155459592847
$ws.Range("J2").Value = "OPTIEMUS"
$ws.Range("K2").Value = -4.393214441061331

# --- Row 3 (Week 3) ---
$ws.Range("C3").Value = 57
$ws.Range("D3").Value = 50
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 87.71929824561403
$ws.Range("G3").Value = 3.66457762960062
$ws.Range("H3").Value = "SWANENERGY"
$ws.Range("I3").Value = 12.82555282555282
$ws.Range("J3").Value = "VISHNU"
$ws.Range("K3").Value = -7.584830339321358

# --- Row 4 (Week 2) ---
$ws.Range("C4").Value = 77
$ws.Range("D4").Value = 64
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 83.11688311688312
$ws.Range("G4").Value = 4.705261222161528
$ws.Range("H4").Value = "SWSOLAR"
$ws.Range("I4").Value = 16.49961449498844
$ws.Range("J4").Value = "JIOFIN"
$ws.Range("K4").Value = -4.272151898734177

# --- Row 5 (Week 1) ---
$ws.Range("C5").Value = 130
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = 30
$ws.Range("F5").Value = 76.92307692307693
$ws.Range("G5").Value = 3.602457862903674
$ws.Range("H5").Value = "MOTHERSON"
$ws.Range("I5").Value = 38.33344090879752
$ws.Range("J5").Value = "ANANDRATHI"
$ws.Range("K5").Value = -22.94392744106948
